# Restore C10 (Integer min for rule R30) from 18 to 1,
# as captured by the target revision's diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
